{"js": "// Update the worksheet date and the twenty-five \"three-digit \u00f7 one-digit\"\n// division problems to the next day's generated set of numbers.\nconst replacements = [\n  [\"2024-02-16 Friday\", \"2024-02-17 Saturday\"],\n  [\"628\u00f72=\", \"981\u00f75=\"],\n  [\"302\u00f79=\", \"340\u00f72=\"],\n  [\"780\u00f74=\", \"289\u00f74=\"],\n  [\"727\u00f79=\", \"339\u00f72=\"],\n  [\"197\u00f76=\", \"942\u00f76=\"],\n  [\"394\u00f79=\", \"503\u00f75=\"],\n  [\"177\u00f73=\", \"972\u00f72=\"],\n  [\"931\u00f75=\", \"672\u00f79=\"],\n  [\"509\u00f77=\", \"862\u00f77=\"],\n  [\"398\u00f73=\", \"793\u00f79=\"],\n  [\"566\u00f74=\", \"717\u00f76=\"],\n  [\"114\u00f78=\", \"113\u00f72=\"],\n  [\"990\u00f74=\", \"633\u00f73=\"],\n  [\"757\u00f74=\", \"936\u00f74=\"],\n  [\"407\u00f73=\", \"418\u00f79=\"],\n  [\"521\u00f74=\", \"768\u00f79=\"],\n  [\"477\u00f73=\", \"391\u00f73=\"],\n  [\"514\u00f74=\", \"277\u00f75=\"],\n  [\"510\u00f73=\", \"848\u00f73=\"],\n  [\"173\u00f76=\", \"889\u00f75=\"],\n  [\"755\u00f73=\", \"574\u00f78=\"],\n  [\"230\u00f74=\", \"668\u00f79=\"],\n  [\"648\u00f77=\", \"552\u00f76=\"],\n  [\"137\u00f74=\", \"816\u00f74=\"],\n  [\"270\u00f77=\", \"229\u00f74=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the twenty-five \"three-digit \u00f7 one-digit\"\n# division problems to the next day's generated set of numbers.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-02-16 Friday\", \"2024-02-17 Saturday\"),\n    @(\"628\u00f72=\", \"981\u00f75=\"),\n    @(\"302\u00f79=\", \"340\u00f72=\"),\n    @(\"780\u00f74=\", \"289\u00f74=\"),\n    @(\"727\u00f79=\", \"339\u00f72=\"),\n    @(\"197\u00f76=\", \"942\u00f76=\"),\n    @(\"394\u00f79=\", \"503\u00f75=\"),\n    @(\"177\u00f73=\", \"972\u00f72=\"),\n    @(\"931\u00f75=\", \"672\u00f79=\"),\n    @(\"509\u00f77=\", \"862\u00f77=\"),\n    @(\"398\u00f73=\", \"793\u00f79=\"),\n    @(\"566\u00f74=\", \"717\u00f76=\"),\n    @(\"114\u00f78=\", \"113\u00f72=\"),\n    @(\"990\u00f74=\", \"633\u00f73=\"),\n    @(\"757\u00f74=\", \"936\u00f74=\"),\n    @(\"407\u00f73=\", \"418\u00f79=\"),\n    @(\"521\u00f74=\", \"768\u00f79=\"),\n    @(\"477\u00f73=\", \"391\u00f73=\"),\n    @(\"514\u00f74=\", \"277\u00f75=\"),\n    @(\"510\u00f73=\", \"848\u00f73=\"),\n    @(\"173\u00f76=\", \"889\u00f75=\"),\n    @(\"755\u00f73=\", \"574\u00f78=\"),\n    @(\"230\u00f74=\", \"668\u00f79=\"),\n    @(\"648\u00f77=\", \"552\u00f76=\"),\n    @(\"137\u00f74=\", \"816\u00f74=\"),\n    @(\"270\u00f77=\", \"229\u00f74=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
